# Chap03: Pictures and caption done.
#
# Applies the "J=3" / "J=2" caption split (-> "M" / "=3" and "M" / "=2")
# and the corresponding text-box widening on slide 1 of Recomb.pptx.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "Text Box 14" carrying the green "J=3" caption (shape #15) ---------
$shpJ3 = $s.Shapes.Item(15)

# Widen the box to make room for the extra character (height/position
# untouched).
$shpJ3.Width = 46.09582710266114

$trJ3 = $shpJ3.TextFrame.TextRange
$trJ3.Text = "M=3"
$mRun = $trJ3.Characters(1, 1)
$mRun.Font.Color.RGB = 5287936        # 00B050
$eqRun = $trJ3.Characters(2, 2)
$eqRun.Font.Color.RGB = 5287936       # 00B050

# --- "Text Box 14" carrying the orange "J=2" caption (shape #37) --------
$shpJ2 = $s.Shapes.Item(37)

# Shift left and widen the box to make room for the extra character.
$shpJ2.Left = 515.2322692871094
$shpJ2.Width = 46.09582710266114

$trJ2 = $shpJ2.TextFrame.TextRange
$trJ2.Text = "M=2"
$mRun2 = $trJ2.Characters(1, 1)
$mRun2.Font.Color.RGB = 49407         # FFC000
$eqRun2 = $trJ2.Characters(2, 2)
$eqRun2.Font.Color.RGB = 49407        # FFC000
